{"js": "const oldText = \"V roku S\u00fahvezdie Bl\u00ed\u017eenci 2022: 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\";\nconst newText = \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Bl\u00ed\u017eenci: 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\";\n\nconst results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$oldText = \"V roku S\u00fahvezdie Bl\u00ed\u017eenci 2022: 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\"\n$newText = \"V roku 2022 m\u00f4\u017eete pozorova\u0165 s\u00fahvezdie S\u00fahvezdie Bl\u00ed\u017eenci: 14. \u2013 23. febru\u00e1ra, 14. \u2013 24. marca\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = $oldText\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = $newText\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n"}
